$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Column C ("Förändrad") holds a date serial number that was bumped by one day
# (45188 -> 45189) for every data row, from row 2 through row 504.
$ws.Range("C2:C504").Value2 = 45189
